# TERRA VERSION, REMOVED RASTER, RGDAL
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# buurtcode value in A2 changed from BU02220303 to BU15810004
$ws.Range("A2").Value = "BU15810004"

# Active cell/selection moved from A3 to A2
$ws.Range("A2").Select()

# Reflect the updated window position recorded by the workbook view
# (best-effort; harmless if the host doesn't persist window geometry)
try {
    $excel.ActiveWindow.Left = 3390
    $excel.ActiveWindow.Top = 630
} catch {
}
